$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.915.13"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "'1.879.70"
$ws.Range("E3").Value = "  -0.37%  "
$ws.Range("E4").Value = "  -0.64%  "
$ws.Range("D5").Value = "'324.86"
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'0.4599"
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("D8").Value = "'0.3882"
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.07854"
$ws.Range("E9").Value = "  -0.30%  "
$ws.Range("D10").Value = "'0.9865"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").Value = "'21.79"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").Value = "'1.831.90"
$ws.Range("E12").Value = "  -4.07%  "
$ws.Range("D13").Value = "'7.002"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "'5.646"
$ws.Range("E14").Value = "  -1.26%  "
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "'88.06"
$ws.Range("E16").Value = "  +0.58%  "
$ws.Range("E17").Value = "  -0.52%  "
$ws.Range("D18").Value = "'0.000009984"
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("E19").Value = "  -1.32%  "
$ws.Range("E20").Value = "  -0.47%  "
$ws.Range("D21").Value = "'28.905.05"
$ws.Range("E21").Value = "  +0.96%  "
$ws.Range("D22").Value = "'5.229"
$ws.Range("E22").Value = "  -2.22%  "
$ws.Range("E23").Value = "  -0.62%  "
$ws.Range("D24").Value = "'2.091"
$ws.Range("E24").Value = "  +1.58%  "
$ws.Range("D25").Value = "'156.39"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "'19.33"
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Value = "'6.019"
$ws.Range("E27").Value = "  +2.49%  "
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("D29").Value = "'117.42"
$ws.Range("E29").Value = "  -0.86%  "
$ws.Range("D30").Value = "'0.09358"
$ws.Range("E30").Value = "  +0.00%  "
$ws.Range("D31").Value = "'0.9021"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("D32").Value = "'5.253"
$ws.Range("E32").Value = "  -1.09%  "
$ws.Range("D33").Value = "'1.313"
$ws.Range("E33").Value = "  -2.10%  "
$ws.Range("D34").Value = "'3.255"
$ws.Range("E34").Value = "  -0.36%  "
$ws.Range("D35").Value = "'1.181"
$ws.Range("E35").Value = "  +1.50%  "
$ws.Range("D36").Value = "'0.05730"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("D37").Value = "'0.02070"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("E38").Value = "  -0.56%  "
$ws.Range("D39").Value = "'7.655"
$ws.Range("E39").Value = "  -6.01%  "
$ws.Range("D40").Value = "'0.5643"
$ws.Range("E40").Value = "  -0.85%  "
$ws.Range("D41").Value = "'0.1764"
$ws.Range("E41").Value = "  -1.93%  "
$ws.Range("D42").Value = "'9.667"
$ws.Range("E42").Value = "  -0.98%  "
$ws.Range("D43").Value = "'2.271"
$ws.Range("E43").Value = "  +4.65%  "
$ws.Range("D44").Value = "'11.90"
$ws.Range("E44").Value = "  +0.74%  "
$ws.Range("D45").Value = "'0.5343"
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "'0.07048"
$ws.Range("E46").Value = "  -1.44%  "
$ws.Range("E47").Value = "  -0.17%  "
$ws.Range("D48").Value = "'2.532"
$ws.Range("E48").Value = "  +1.67%  "
$ws.Range("D49").Value = "'112.49"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'1.061"
$ws.Range("E50").Value = "  -5.18%  "
$ws.Range("D51").Value = "'70.56"
$ws.Range("E51").Value = "  -0.60%  "
